# Add two new trailing rows (2021年, 2022年) to Sheet1, mirroring the
# existing table's layout:
#   A = year label, B = 住房保养维修及管理类, C = 居住类,
#   D = 水电燃料类, E = 租赁房房租类, F = 自有住房类
#
# Row 8 (2022年) only has data published for column C so far; columns
# B, D, E, F are "present but blank" cells (same convention already used
# by F5/F6 in the original sheet - an empty text cell, not simply an
# absent one). Assigning a bare "" via COM clears/removes the cell
# instead of leaving an empty string behind, so we use Excel's classic
# leading single-quote ("text prefix") trick - entering just a quote
# character commits an empty *text* cell. We then paste only the number
# formatting (not the value) from a plain, unstyled cell on top, so the
# quote-prefix styling doesn't linger on the cell itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 2021年 ---------------------------------------------------
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 102.3
$ws.Range("C7").Value = 100.8
$ws.Range("D7").Value = 101.6
$ws.Range("E7").Value = 100.4
$ws.Range("F7").Value = "'"   # blank/empty text cell

# --- Row 8: 2022年 ---------------------------------------------------
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = "'"   # blank/empty text cell
$ws.Range("C8").Value = 100.5
$ws.Range("D8").Value = "'"   # blank/empty text cell
$ws.Range("E8").Value = "'"   # blank/empty text cell
$ws.Range("F8").Value = "'"   # blank/empty text cell

# Give the new year labels (column A) the same bold/centered/bordered
# style already used by the other year cells (A2:A6).
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)   # xlPasteFormats

# Strip the quote-prefix styling picked up from the "'" blank-cell
# trick above, by pasting the plain (unstyled) format from an existing
# data cell over the top - this leaves the values untouched.
$ws.Range("B6").Copy()
$ws.Range("F7").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("B8").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("D8").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("E8").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("F8").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = $false
